$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the third row entirely (A3, B3) and the extra C2 cell
$ws.Range("C2").Value = $null
$ws.Range("A3:B3").Value = $null

# Update B2 value to the new course/roll number
$ws.Range("B2").Value = "B20EE011"
